# ============================================================================
# Reproduces the diff:
#   1. Splits "...Такая графика вам поможет погрузиться в " into several
#      runs ("Такая графика " / "поможет вам" / " погрузиться " /
#      "в мир игр."), adds a line break, and inserts the picture
#      (docPr id="2" name="Рисунок 2") right after that paragraph -- i.e.
#      the picture that used to sit at the very end of the document moves
#      up here.
#   2. Inserts a lastRenderedPageBreak + two line breaks right before the
#      "7 Платформа" heading run.
#   3. Removes the picture (docPr id="1" name="Рисунок 1") that used to be
#      in the final paragraph of the document, leaving that paragraph
#      holding the relocated "_GoBack" bookmark instead.
#
# Bookmark w:id numbers are NOT hard-coded to the final values seen in the
# diff (0 / 1): the engine renumbers bookmark ids sequentially by document
# order whenever the package is (re)written, exactly like real Word does on
# save, so only the *locations* of the bookmarks matter here.
# ============================================================================

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
            'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
            'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
            'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' +
            'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' +
            'xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing">' +
        '<w:body>' + $bodyInner + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Step 0: drop the "_GoBack" bookmark from its old location (right after
# "...графика вам поможет погрузиться в "); it gets re-created at the new
# end-of-document location in step 3.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# Step 1: rewrite "...Такая графика вам поможет погрузиться в " into the
# multi-run form and insert the relocated picture + line break after it.
# ---------------------------------------------------------------------
$target1 = "Такая графика вам поможет погрузиться в "
$search1 = $d.Content
$found1 = $search1.Find.Execute($target1)
if (-not $found1) { throw "Could not find target1 text" }
# Build a *fresh* Range from the found bounds -- reusing the Range that
# Find itself mutated as the InsertXML receiver swallows the whole host
# paragraph instead of just the matched span.
$insRange1 = $d.Range($search1.Start, $search1.End)

$body1 = '<w:p>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Такая графика </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>поможет вам</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> погрузиться </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>в мир игр.</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:br/></w:r>' +
  '<w:r><w:rPr><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
    '<w:drawing>' +
      '<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="3888391F" wp14:editId="5D042EB4">' +
        '<wp:extent cx="5935980" cy="3337560"/>' +
        '<wp:effectExtent l="0" t="0" r="7620" b="0"/>' +
        '<wp:docPr id="2" name="&#1056;&#1080;&#1089;&#1091;&#1085;&#1086;&#1082; 2"/>' +
        '<wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr>' +
        '<a:graphic>' +
          '<a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
            '<pic:pic>' +
              '<pic:nvPicPr>' +
                '<pic:cNvPr id="0" name="Picture 1"/>' +
                '<pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr>' +
              '</pic:nvPicPr>' +
              '<pic:blipFill>' +
                '<a:blip r:embed="rId5">' +
                  '<a:extLst>' +
                    '<a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext>' +
                  '</a:extLst>' +
                '</a:blip>' +
                '<a:srcRect/>' +
                '<a:stretch><a:fillRect/></a:stretch>' +
              '</pic:blipFill>' +
              '<pic:spPr bwMode="auto">' +
                '<a:xfrm><a:off x="0" y="0"/><a:ext cx="5935980" cy="3337560"/></a:xfrm>' +
                '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom>' +
                '<a:noFill/>' +
                '<a:ln><a:noFill/></a:ln>' +
              '</pic:spPr>' +
            '</pic:pic>' +
          '</a:graphicData>' +
        '</a:graphic>' +
      '</wp:inline>' +
    '</w:drawing>' +
  '</w:r>' +
  '</w:p>'

$insRange1.InsertXML( (New-PkgXml $body1) )

# ---------------------------------------------------------------------
# Step 2: insert lastRenderedPageBreak + two breaks before "7 Платформа"
# ---------------------------------------------------------------------
$search2 = $d.Content
$found2 = $search2.Find.Execute("7 Платформа")
if (-not $found2) { throw "Could not find '7 Платформа' heading" }
$insPoint2 = $d.Range($search2.Start, $search2.Start)

$body2 = '<w:p>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:lastRenderedPageBreak/><w:br/></w:r>' +
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:br/></w:r>' +
  '</w:p>'

$insPoint2.InsertXML( (New-PkgXml $body2) )

# ---------------------------------------------------------------------
# Step 3: remove the picture (docPr id="1", "Рисунок 1") from the final
# paragraph of the document and drop the relocated "_GoBack" bookmark in
# its place.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$delRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$delRange.Text = ""

$lastPara2 = $d.Paragraphs.Last
$lastRange2 = $lastPara2.Range
$insPoint3 = $d.Range($lastRange2.Start, $lastRange2.Start)
$body3 = '<w:p><w:bookmarkStart w:id="5" w:name="_GoBack"/><w:bookmarkEnd w:id="5"/></w:p>'
$insPoint3.InsertXML( (New-PkgXml $body3) )

Write-Output "All edits applied"
